$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text, matching the
# original inline-string cell content (column D holds text prices).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.263.69"
$ws.Range("E2").Value = "  +6.59%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.009.19"
$ws.Range("E3").Value = "  +3.56%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.92"
$ws.Range("E5").Value = "  +2.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.10"
$ws.Range("E6").Value = "  +13.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.005.30"
$ws.Range("E8").Value = "  +3.51%  "

$ws.Range("E9").Value = "  +3.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("E10").Value = "  -4.53%  "

$ws.Range("E11").Value = "  +4.06%  "

$ws.Range("E12").Value = "  +5.33%  "

$ws.Range("E13").Value = "  +6.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.58"
$ws.Range("E14").Value = "  +6.13%  "

$ws.Range("E15").Value = "  -0.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.219.13"
$ws.Range("E16").Value = "  +6.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.508.31"
$ws.Range("E17").Value = "  +3.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.92"
$ws.Range("E18").Value = "  +5.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.006.96"
$ws.Range("E19").Value = "  +3.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "454.77"
$ws.Range("E20").Value = "  +5.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.86"
$ws.Range("E21").Value = "  +6.29%  "

$ws.Range("E22").Value = "  +4.54%  "

$ws.Range("E23").Value = "  +7.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.21"
$ws.Range("E24").Value = "  +4.61%  "

$ws.Range("E25").Value = "  +14.54%  "

$ws.Range("E26").Value = "  +2.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.53"
$ws.Range("E27").Value = "  +4.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.14"
$ws.Range("E29").Value = "  +17.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.40"
$ws.Range("E30").Value = "  +19.45%  "

$ws.Range("E31").Value = "  -6.09%  "

$ws.Range("E32").Value = "  +4.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.21"
$ws.Range("E33").Value = "  +5.98%  "

$ws.Range("E34").Value = "  +4.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.991"
$ws.Range("E36").Value = "  +3.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.79"
$ws.Range("E37").Value = "  +7.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.18"
$ws.Range("E38").Value = "  +15.44%  "

$ws.Range("E39").Value = "  +1.94%  "

$ws.Range("E40").Value = "  +2.20%  "

$ws.Range("E41").Value = "  +16.06%  "

$ws.Range("E42").Value = "  +8.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "44.03"
$ws.Range("E43").Value = "  +8.05%  "

$ws.Range("E44").Value = "  +3.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "394.52"
$ws.Range("E45").Value = "  +13.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0359"
$ws.Range("E46").Value = "  +6.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.791.87"
$ws.Range("E47").Value = "  +2.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.85"
$ws.Range("E48").Value = "  +0.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.94"
$ws.Range("E50").Value = "  +11.84%  "

$ws.Range("E51").Value = "  +4.34%  "
